$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row4 = @(0.5867794862083637, 0.588, 0.5911496908674081, 0.59, 0.6755282899659024, 0.6990000000000001, 0.6588925889241295, 0.664, 0.4888536848760653, 0.4970000000000001, 0.4853343332709336, 0.483, 0.6580453174800088, 0.6859999999999999, 0.6385261638611831, 0.647, 0.5910910638846751, 0.5940000000000001, 0.5930100494528215, 0.5924999999999999, 0.6714857468567652, 0.696, 0.654044396947697, 0.6599999999999999, 0.6623944431189186, 0.6910000000000001, 0.6416551457567504, 0.6505)
for ($i = 0; $i -lt $row4.Length; $i++) {
  $ws.Cells.Item(4, 2 + $i).Value = $row4[$i]
}

$row5 = @(0.6524668262824616, 0.8470000000000001, 0.5342609065801748, 0.5525, 0.6623177022921712, 0.849, 0.5511598390280821, 0.5740000000000001, 0.6687938021972293, 0.946, 0.5182224495543275, 0.532, 0.6540053142392115, 0.852, 0.5355237034507258, 0.5534999999999999, 0.6561070483877145, 0.86, 0.5339164801115952, 0.553, 0.6702176239762629, 0.885, 0.5450095930567206, 0.57, 0.6553253939830184, 0.861, 0.5334430843523295, 0.5515000000000001)
for ($i = 0; $i -lt $row5.Length; $i++) {
  $ws.Cells.Item(5, 2 + $i).Value = $row5[$i]
}

$row6 = @(0.5743952573160535, 0.562, 0.5997305832069287, 0.6010000000000001, 0.6935710155645171, 0.6960000000000001, 0.7070314191332031, 0.6945, 0.5193020752264974, 0.5240000000000001, 0.5244858522170122, 0.5225000000000001, 0.6676778392119738, 0.6499999999999999, 0.710395504594944, 0.69, 0.5829022968259514, 0.5700000000000001, 0.6113618683618081, 0.6105, 0.6853297188491376, 0.6969999999999998, 0.6863753167185024, 0.683, 0.6795534778210511, 0.6599999999999999, 0.7226943761333404, 0.7)
for ($i = 0; $i -lt $row6.Length; $i++) {
  $ws.Cells.Item(6, 2 + $i).Value = $row6[$i]
}

